# "Generate Report for Handoff" - update the localization-status report with
# the latest handoff timestamps for the most recently processed file
# (6e85a030-a4e5-4092-8df3-fdfbde5c9ecd.md).

$wb = $excel.ActiveWorkbook

# zh-cn: Latest Handoff Datetime (column H) for row 5
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2016-11-14 05:58:51"

# de-de: Latest Handoff Datetime (column H) for row 5
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2016-11-14 05:59:04"

# Overview: Latest HO Xliff Generate Date (column G) for row 5 -
# reflects the newest of the per-language handoff datetimes above.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2016-11-14 05:59:04"
